$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply refreshed price/volume data (and a 3-way coin reorder for rows 45/46/49).
# Price-column values that look numeric are written via a text number format
# so Excel keeps exact digits (e.g. trailing zeros) instead of coercing to a
# float, then the style is reset to Normal to match the original unstyled cells.

# Row 2
$ws.Range("D2").Value = "29.372.69"
$ws.Range("E2").Value = "  -0.13%  "

# Row 3
$ws.Range("D3").Value = "1.849.11"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6282"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07612"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2915"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.68%  "

# Row 10
$ws.Range("E10").Value = "  +0.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07749"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("E12").Value = "  -0.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6805"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001050"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.129"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "29.426.49"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("E19").Value = "  -1.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.470"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "158.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1390"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.445"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.448"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.84%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.474"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "

# Row 29
$ws.Range("E29").Value = "  -2.04%  "

# Row 30
$ws.Range("E30").Value = "  -0.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.062"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.832"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.586"
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "1.235.81"
$ws.Range("E36").Value = "  +1.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01807"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.727"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.429"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9053"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "

# Row 44
$ws.Range("E44").Value = "  +0.76%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.42%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.045"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1156"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.683"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000115"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05700"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4628"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
